$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new column D ("Assembly") - shifts existing D..I to E..J ---
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).ColumnWidth = 35

# --- Header row ---
$ws.Range("D5").Value = "Assembly"

# --- Existing rows: fill in the new Assembly column ---
$ws.Range("D7").Value  = "External case cables & Lower Printer Case"
$ws.Range("D9").Value  = "External case cables & Lower Printer Case"
$ws.Range("D11").Value = "External case cables & Lower Printer Case"
$ws.Range("D13").Value = "External case cables & Lower Printer Case"
$ws.Range("D15").Value = "External case cables & Lower Printer Case"
$ws.Range("D17").Value = "External case cables & Lower Printer Case"

$ws.Range("D50").Value = "General"
$ws.Range("D52").Value = "General"
$ws.Range("D54").Value = "General"

# --- New BOM rows ---

# Row 57-58: item 25, PWR ENT MOD RCPT IEC320-C14 PNL
$ws.Range("B57").Value = 25
$ws.Range("C57").Value = "PWR ENT MOD RCPT IEC320-C14 PNL"
$ws.Range("D57").Value = "External Case Back Cover"
$ws.Range("E57").Value = 1
$ws.Range("F57").Value = "Digikey"
$ws.Range("G57").Value = "Q307-ND"
$ws.Range("H57").Value = "Short"
$ws.Range("C58").Value = "10A"

# Row 60-61: item 26, FUSE GLASS 5A 250VAC 5X20MM
$ws.Range("B60").Value = 26
$ws.Range("C60").Value = "FUSE GLASS 5A 250VAC 5X20MM"
$ws.Range("D60").Value = "External Case Back Cover"
$ws.Range("E60").Value = 2
$ws.Range("F60").Value = "Digikey"
$ws.Range("G60").Value = "F1748-ND"
$ws.Range("H60").Value = "Short"
$ws.Range("C61").Value = "(or similar 5A of this size)"
$ws.Range("E61").Value = "Get spares"
$ws.Range("E61").HorizontalAlignment = -4108

# Row 63: item 24, SWITCH PUSHBUTTON SPST 5A 125V
$ws.Range("B63").Value = 24
$ws.Range("C63").Value = "SWITCH PUSHBUTTON SPST 5A 125V"
$ws.Range("D63").Value = "External Case Front Cover"
$ws.Range("E63").Value = 1
$ws.Range("F63").Value = "Digikey"
$ws.Range("H63").Value = "Short"

# Row 65: item 26, 24 VDC Fan - looking for specs
$ws.Range("B65").Value = 26
$ws.Range("C65").Value = "24 VDC Fan - looking for specs"
$ws.Range("D65").Value = "External Case Front Cover"
$ws.Range("E65").Value = 1
$ws.Range("F65").Value = "Digikey"
$ws.Range("G65").Value = "?"
$ws.Range("H65").Value = "Short"

# --- Update view selection ---
$ws.Range("D67").Select()
